$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Snapshot the source rows (columns A to AY) before any writes, since the row
#     permutation below forms a single long cycle and every row is both a source
#     and a destination. ---
$rowData = @{}
$rowData[2] = $ws.Range("A2:AY2").Value()
$rowData[3] = $ws.Range("A3:AY3").Value()
$rowData[4] = $ws.Range("A4:AY4").Value()
$rowData[5] = $ws.Range("A5:AY5").Value()
$rowData[6] = $ws.Range("A6:AY6").Value()
$rowData[7] = $ws.Range("A7:AY7").Value()
$rowData[9] = $ws.Range("A9:AY9").Value()
$rowData[10] = $ws.Range("A10:AY10").Value()
$rowData[11] = $ws.Range("A11:AY11").Value()
$rowData[12] = $ws.Range("A12:AY12").Value()
$rowData[13] = $ws.Range("A13:AY13").Value()
$rowData[14] = $ws.Range("A14:AY14").Value()
$rowData[15] = $ws.Range("A15:AY15").Value()
$rowData[16] = $ws.Range("A16:AY16").Value()
$rowData[17] = $ws.Range("A17:AY17").Value()
$rowData[18] = $ws.Range("A18:AY18").Value()
$rowData[19] = $ws.Range("A19:AY19").Value()
$rowData[20] = $ws.Range("A20:AY20").Value()
$rowData[21] = $ws.Range("A21:AY21").Value()
$rowData[22] = $ws.Range("A22:AY22").Value()
$rowData[23] = $ws.Range("A23:AY23").Value()
$rowData[24] = $ws.Range("A24:AY24").Value()
$rowData[26] = $ws.Range("A26:AY26").Value()
$rowData[27] = $ws.Range("A27:AY27").Value()
$rowData[28] = $ws.Range("A28:AY28").Value()
$rowData[29] = $ws.Range("A29:AY29").Value()
$rowData[30] = $ws.Range("A30:AY30").Value()
$rowData[31] = $ws.Range("A31:AY31").Value()
$rowData[32] = $ws.Range("A32:AY32").Value()
$rowData[33] = $ws.Range("A33:AY33").Value()

# --- Columns whose text content could be auto-converted by Excel (dates/times/
#     numeric-looking text) if written without first forcing a Text number format. ---
$textCols = @("I", "Y", "Z", "AA", "AB")

function Set-RowFromSnapshot($destRow, $srcRow) {
    foreach ($col in $textCols) {
        $ws.Range("$col$destRow").NumberFormat = "@"
    }
    $ws.Range("A$destRow`:AY$destRow").Value() = $rowData[$srcRow]
}

# --- Apply the row permutation described by the diff (rows 8 and 25 are left untouched) ---
Set-RowFromSnapshot 2 27
Set-RowFromSnapshot 3 33
Set-RowFromSnapshot 4 31
Set-RowFromSnapshot 5 29
Set-RowFromSnapshot 6 15
Set-RowFromSnapshot 7 21
Set-RowFromSnapshot 9 10
Set-RowFromSnapshot 10 11
Set-RowFromSnapshot 11 24
Set-RowFromSnapshot 12 2
Set-RowFromSnapshot 13 19
Set-RowFromSnapshot 14 23
Set-RowFromSnapshot 15 16
Set-RowFromSnapshot 16 22
Set-RowFromSnapshot 17 5
Set-RowFromSnapshot 18 32
Set-RowFromSnapshot 19 6
Set-RowFromSnapshot 20 14
Set-RowFromSnapshot 21 17
Set-RowFromSnapshot 22 3
Set-RowFromSnapshot 23 13
Set-RowFromSnapshot 24 30
Set-RowFromSnapshot 26 7
Set-RowFromSnapshot 27 28
Set-RowFromSnapshot 28 9
Set-RowFromSnapshot 29 18
Set-RowFromSnapshot 30 26
Set-RowFromSnapshot 31 12
Set-RowFromSnapshot 32 20
Set-RowFromSnapshot 33 4

Write-Output "Row permutation applied."
